$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New meeting minutes row (row 11), entered date/time first, as the author did
$ws.Range("A10").Copy($ws.Range("A11"))
$ws.Range("A11").Value = "October 21 2023"

$ws.Range("C10").Copy($ws.Range("C11"))
$ws.Range("C11").Value = "7:10PM"

# Fill in the remainder of row 10 (previously only had A10:C10 populated)
$ws.Range("C10").Copy($ws.Range("D10"))
$ws.Range("D10").Value = "8:00PM"

$ws.Range("A10").Copy($ws.Range("E10"))
$ws.Range("E10").Value = "Finish up Presentation Slides"

$ws.Range("A10").Copy($ws.Range("E11"))
$ws.Range("E11").Value = "Figuring out what are we supposed to do next in the project"

$ws.Range("C10").Copy($ws.Range("D11"))
$ws.Range("D11").Value = "8:35PM"

$ws.Range("B10").Copy($ws.Range("B11"))
$ws.Range("B11").Value = "Sedat, Madison, Joseph, David, Sean"

# Update selection to match the new last-edited cell
$ws.Range("D11").Select()

# Nudge the workbook window position, matching the saved view state
$wb.Windows.Item(1).Left = 2400
$wb.Windows.Item(1).Top = 1080
